$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Paragraph 4 ("This research aims to enhance ...") - rewrite the
#    tail of the paragraph: "adapted obstacle designs. The questionnaire
#    design includes ..." becomes "dynamic obstacle designs. The
#    questionnaire was designed using UEQ, ...". Re-set the whole
#    paragraph's text (minus its trailing paragraph mark) so Word
#    regenerates a single clean run (and drops the old spell-check
#    proofErr markers around "MoodME").
# ------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.End = $r4.End - 1
$r4.Text = "This research aims to enhance the quality of horror game players' experience by adjusting the obstacles present in the game levels. The adjustments are made by analyzing the facial emotional data of players during gaming sessions. The methodology employed in this research is the waterfall model, where all processes must be carried out sequentially, starting from the basic game creation to processing data from the MoodME library, followed by White Box and Black Box testing. The game is evaluated using a questionnaire given to players who have experienced the game with dynamic obstacle designs. The questionnaire was designed using UEQ, which focuses on measuring the user experience needed to gauge how well the game can adapt obstacles to player performance and the attractiveness of the game design."

# ------------------------------------------------------------------
# 2. Paragraph 5 ("Based on questionnaire data from 27 individuals ...")
#    gets entirely new closing content about the updated survey numbers.
# ------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$r5 = $p5.Range
$r5.End = $r5.End - 1
$r5.Text = "Based on questionnaire data from 30 individuals, it shows that the use of obstacle adjustment concepts that can adapt difficulty levels based on players' facial expressions has successfully attracted players repeatedly (mean attractiveness = 1.59). However, there is a need to provide further explanation about the purpose and function of obstacle adjustments based on player emotions (mean clarity = 1.14), as well as ensuring accuracy (mean precision = 1.32) and efficiency (mean efficiency = 1.47) in adjusting difficulty levels to player performance. The use of obstacles that can adapt to players has successfully stimulated (mean stimulation = 1.71) players to continue playing until completion. The feature of adjusting obstacles based on player scores and emotions is considered innovative and up-to-date (mean novelty = 1.45). Thus, this game has great potential to maintain player interest with high attractiveness, provided special attention is given to the development of specific aspects that affect player experience."

# ------------------------------------------------------------------
# 3. Remove the trailing empty paragraph that used to sit between the
#    abstract text and the section break.
# ------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$p6 = $d.Paragraphs(6)
$delRange = $d.Range($p5.Range.End - 1, $p6.Range.End)
$delRange.Delete()

# ------------------------------------------------------------------
# 4. Bump the starting page number (lower-roman) for the section from
#    4 to 5.
# ------------------------------------------------------------------
$sec = $d.Sections(1)
$sec.Footers(1).PageNumbers.StartingNumber = 5
